# Weekly update: a new Fruta/Palta price record for Agrícola del Norte S.A.
# de Arica needs to be inserted as row 21 (pushing the existing row 21..85
# records down to 22..86, dimension growing from A1:T85 to A1:T86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 21..85 down one row, leaving a blank (but format-inheriting) row 21.
$ws.Rows(21).Insert()

# Populate the newly inserted row 21 with the new observation.
$ws.Cells.Item(21, 1).Value2  = 1
$ws.Cells.Item(21, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value2  = 44659
$ws.Cells.Item(21, 5).Value2  = 15
$ws.Cells.Item(21, 6).Value2  = "Fruta"
$ws.Cells.Item(21, 7).Value2  = 100106
$ws.Cells.Item(21, 8).Value2  = "Oleaginosos"
$ws.Cells.Item(21, 9).Value2  = 100106002
$ws.Cells.Item(21, 10).Value2 = "Palta"
$ws.Cells.Item(21, 11).Value2 = "Hass"
$ws.Cells.Item(21, 12).Value2 = "Primera"
$ws.Cells.Item(21, 13).Value2 = 200
$ws.Cells.Item(21, 14).Value2 = 25000
$ws.Cells.Item(21, 15).Value2 = 27000
$ws.Cells.Item(21, 16).Value2 = 26000
$ws.Cells.Item(21, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(21, 18).Value2 = "Perú"
$ws.Cells.Item(21, 19).Value2 = 2600
$ws.Cells.Item(21, 20).Value2 = 10
